# Update cryptocurrency price (D) and 1h volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.716.31"
$ws.Range("E2").Value = "'  +0.38%  "
$ws.Range("D3").Value = "'1.852.25"
$ws.Range("E3").Value = "'  +0.56%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "'  +0.21%  "
$ws.Range("D5").Value = "'312.61"
$ws.Range("E5").Value = "'  -0.58%  "
$ws.Range("E6").Value = "'  +0.19%  "
$ws.Range("D7").Value = "'0.4285"
$ws.Range("E7").Value = "'  +1.52%  "
$ws.Range("D8").Value = "'0.3588"
$ws.Range("E8").Value = "'  -1.32%  "
$ws.Range("D9").Value = "'0.07299"
$ws.Range("E9").Value = "'  +0.12%  "
$ws.Range("D10").Value = "'0.8759"
$ws.Range("E10").Value = "'  -1.41%  "
$ws.Range("E11").Value = "'  +0.55%  "
$ws.Range("D12").Value = "'1.870.73"
$ws.Range("E12").Value = "'  +3.25%  "
$ws.Range("D13").Value = "'6.552"
$ws.Range("E13").Value = "'  -0.12%  "
$ws.Range("D14").Value = "'5.336"
$ws.Range("E14").Value = "'  +0.09%  "
$ws.Range("D15").Value = "'0.07001"
$ws.Range("E15").Value = "'  +1.60%  "
$ws.Range("D16").Value = "'1.007"
$ws.Range("E16").Value = "'  +0.38%  "
$ws.Range("D17").Value = "'79.63"
$ws.Range("E17").Value = "'  +0.59%  "
$ws.Range("D18").Value = "'0.000008963"
$ws.Range("E18").Value = "'  +0.83%  "
$ws.Range("D19").Value = "'1.004"
$ws.Range("E19").Value = "'  +0.39%  "
$ws.Range("E20").Value = "'  -0.83%  "
$ws.Range("D21").Value = "'27.684.01"
$ws.Range("E21").Value = "'  +0.28%  "
$ws.Range("D22").Value = "'5.004"
$ws.Range("E22").Value = "'  +0.43%  "
$ws.Range("D23").Value = "'10.41"
$ws.Range("E23").Value = "'  -1.54%  "
$ws.Range("D24").Value = "'2.114.55"
$ws.Range("E24").Value = "'  +3.59%  "
$ws.Range("D25").Value = "'1.991"
$ws.Range("E25").Value = "'  +3.85%  "
$ws.Range("D26").Value = "'155.18"
$ws.Range("E26").Value = "'  +0.88%  "
$ws.Range("D27").Value = "'18.50"
$ws.Range("E27").Value = "'  -3.07%  "
$ws.Range("D28").Value = "'120.56"
$ws.Range("E28").Value = "'  -1.71%  "
$ws.Range("D29").Value = "'5.283"
$ws.Range("E29").Value = "'  +0.04%  "
$ws.Range("D30").Value = "'1.890"
$ws.Range("E30").Value = "'  -0.01%  "
$ws.Range("D31").Value = "'0.08918"
$ws.Range("E31").Value = "'  -0.16%  "
$ws.Range("D32").Value = "'0.7588"
$ws.Range("E32").Value = "'  -1.16%  "
$ws.Range("D33").Value = "'2.973"
$ws.Range("E33").Value = "'  +1.48%  "
$ws.Range("D34").Value = "'4.519"
$ws.Range("E34").Value = "'  -1.15%  "
$ws.Range("D35").Value = "'1.124"
$ws.Range("E35").Value = "'  +2.66%  "
$ws.Range("D37").Value = "'0.05445"
$ws.Range("E37").Value = "'  +1.36%  "
$ws.Range("D38").Value = "'1.106"
$ws.Range("E38").Value = "'  +0.36%  "
$ws.Range("D39").Value = "'0.01934"
$ws.Range("E39").Value = "'  -0.36%  "
$ws.Range("D40").Value = "'2.831"
$ws.Range("E40").Value = "'  -0.01%  "
$ws.Range("D41").Value = "'0.1668"
$ws.Range("E41").Value = "'  +0.92%  "
$ws.Range("D42").Value = "'0.5089"
$ws.Range("E42").Value = "'  -0.02%  "
$ws.Range("D43").Value = "'6.628"
$ws.Range("E43").Value = "'  -3.50%  "
$ws.Range("D44").Value = "'8.414"
$ws.Range("E44").Value = "'  +1.83%  "
$ws.Range("D45").Value = "'106.40"
$ws.Range("E45").Value = "'  +2.13%  "
$ws.Range("D46").Value = "'0.06531"
$ws.Range("E46").Value = "'  -1.03%  "
$ws.Range("D47").Value = "'10.37"
$ws.Range("E47").Value = "'  -0.56%  "
$ws.Range("E48").Value = "'  -1.11%  "
$ws.Range("D49").Value = "'1.003"
$ws.Range("E49").Value = "'  +0.23%  "
$ws.Range("D50").Value = "'1.626"
$ws.Range("E50").Value = "'  -0.41%  "
$ws.Range("E51").Value = "'  +2.17%  "
